# Standardise "cost_variable" -> "cost_variable_om" in column C (rows 10-39)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C10:C39").Value = "cost_variable_om"

# Update the active selection to match the edited range, as Excel would
# after selecting/editing these cells.
$ws.Range("C10:C39").Select()
